$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1564.8889
$ws.Range("I88").Value = 1400
$ws.Range("J88").Value = 1585.5
$ws.Range("K88").Value = 1400
$ws.Range("L88").Value = 1585.5
$ws.Range("M88").Value = -994
$ws.Range("N88").Value = -2397.5
$ws.Range("H91").Value = 1564.8889
$ws.Range("I91").Value = 1400
$ws.Range("J91").Value = 1585.5
$ws.Range("K91").Value = 1400
$ws.Range("L91").Value = 1585.5
$ws.Range("M91").Value = 4
$ws.Range("N91").Value = -4393.5
$ws.Range("H121").Value = 1226.5385
$ws.Range("I121").Value = 785
$ws.Range("J121").Value = 1263.3334
$ws.Range("K121").Value = 2355
$ws.Range("L121").Value = 3790.0002
$ws.Range("M121").Value = -608
$ws.Range("N121").Value = -7284.0002
$ws.Range("H131").Value = 1152.6666
$ws.Range("I131").Value = 724.1667
$ws.Range("K131").Value = 2172.5001
$ws.Range("M131").Value = 2867.4999
$ws.Range("H138").Value = 2355.6338
$ws.Range("I138").Value = 1183.8572
$ws.Range("J138").Value = 4965.5
$ws.Range("K138").Value = 3551.5716
$ws.Range("L138").Value = 14896.5
$ws.Range("M138").Value = 1588.4284
$ws.Range("N138").Value = -25176.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 972.6667
$ws.Range("I2").Value = 919.1111
$ws.Range("K2").Value = 919.1111
$ws.Range("M2").Value = -806.1111
$ws.Range("H61").Value = 1062
$ws.Range("I61").Value = 852.0625
$ws.Range("J61").Value = 2405.6
$ws.Range("K61").Value = 852.0625
$ws.Range("L61").Value = 2405.6
$ws.Range("M61").Value = -640.0625
$ws.Range("N61").Value = -2829.6
$ws.Range("H74").Value = 1122.0698
$ws.Range("I74").Value = 606.1842
$ws.Range("J74").Value = 5042.8
$ws.Range("K74").Value = 606.1842
$ws.Range("L74").Value = 5042.8
$ws.Range("M74").Value = 267.8158
$ws.Range("N74").Value = -6790.8
$ws.Range("H77").Value = 1122.0698
$ws.Range("I77").Value = 606.1842
$ws.Range("J77").Value = 5042.8
$ws.Range("K77").Value = 3030.921
$ws.Range("L77").Value = 25214
$ws.Range("M77").Value = 1337.079
$ws.Range("N77").Value = -33950
$ws.Range("H116").Value = 972.6667
$ws.Range("I116").Value = 919.1111
$ws.Range("K116").Value = 919.1111
$ws.Range("M116").Value = 1374.8889
$ws.Range("H132").Value = 1699.8545
$ws.Range("I132").Value = 1607.6364
$ws.Range("J132").Value = 2068.7273
$ws.Range("K132").Value = 4822.9092
$ws.Range("L132").Value = 6206.1819
$ws.Range("M132").Value = -2292.9092
$ws.Range("N132").Value = -11266.1819
$ws.Range("H136").Value = 1062
$ws.Range("I136").Value = 852.0625
$ws.Range("J136").Value = 2405.6
$ws.Range("K136").Value = 2556.1875
$ws.Range("L136").Value = 7216.799999999999
$ws.Range("M136").Value = -6.1875
$ws.Range("N136").Value = -12316.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 972.6667
$ws.Range("I3").Value = 919.1111
$ws.Range("K3").Value = 919.1111
$ws.Range("M3").Value = -805.1111
$ws.Range("H68").Value = 30295
$ws.Range("J68").Value = 30295
$ws.Range("L68").Value = 30295
$ws.Range("N68").Value = -31917
$ws.Range("H71").Value = 30295
$ws.Range("J71").Value = 30295
$ws.Range("L71").Value = 90885
$ws.Range("N71").Value = -98997
$ws.Range("H105").Value = 3898.9
$ws.Range("I105").Value = 3887.6667
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 3887.6667
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -2140.6667
$ws.Range("N105").Value = -7494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19538.104
$ws.Range("I31").Value = 24777.523
$ws.Range("J31").Value = 3071.3572
$ws.Range("K31").Value = 24777.523
$ws.Range("L31").Value = 3071.3572
$ws.Range("M31").Value = -24482.523
$ws.Range("N31").Value = -3661.3572
$ws.Range("H34").Value = 19538.104
$ws.Range("I34").Value = 24777.523
$ws.Range("J34").Value = 3071.3572
$ws.Range("K34").Value = 24777.523
$ws.Range("L34").Value = 3071.3572
$ws.Range("M34").Value = -24575.523
$ws.Range("N34").Value = -3475.3572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1362.3889
$ws.Range("I5").Value = 327.625
$ws.Range("J5").Value = 2190.2
$ws.Range("K5").Value = 982.875
$ws.Range("L5").Value = 6570.599999999999
$ws.Range("M5").Value = -870.875
$ws.Range("N5").Value = -6794.599999999999
$ws.Range("H113").Value = 903.9697
$ws.Range("I113").Value = 1209.9474
$ws.Range("J113").Value = 488.7143
$ws.Range("K113").Value = 3629.8422
$ws.Range("L113").Value = 1466.1429
$ws.Range("M113").Value = -1459.8422
$ws.Range("N113").Value = -5806.1429
$ws.Range("H122").Value = 1094.1765
$ws.Range("I122").Value = 818.8333
$ws.Range("J122").Value = 1755
$ws.Range("K122").Value = 7369.4997
$ws.Range("L122").Value = 15795
$ws.Range("M122").Value = -4919.4997
$ws.Range("N122").Value = -20695
$ws.Range("H135").Value = 1362.3889
$ws.Range("I135").Value = 327.625
$ws.Range("J135").Value = 2190.2
$ws.Range("K135").Value = 2948.625
$ws.Range("L135").Value = 19711.8
$ws.Range("M135").Value = -413.625
$ws.Range("N135").Value = -24781.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31748
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -98736
$ws.Range("H139").Value = 49733.332
$ws.Range("J139").Value = 49733.332
$ws.Range("L139").Value = 49733.332
$ws.Range("N139").Value = -60013.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1994.5209
$ws.Range("I132").Value = 1893.075
$ws.Range("J132").Value = 2501.75
$ws.Range("K132").Value = 5679.225
$ws.Range("L132").Value = 7505.25
$ws.Range("M132").Value = -3149.225
$ws.Range("N132").Value = -12565.25
$ws.Range("H133").Value = 25342.46
$ws.Range("J133").Value = 25342.46
$ws.Range("L133").Value = 25342.46
$ws.Range("N133").Value = -30402.46
$ws.Range("H136").Value = 2298.8306
$ws.Range("I136").Value = 1666.186
$ws.Range("J136").Value = 3999.0625
$ws.Range("K136").Value = 4998.558
$ws.Range("L136").Value = 11997.1875
$ws.Range("M136").Value = -2448.558
$ws.Range("N136").Value = -17097.1875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 54178.625
$ws.Range("I46").Value = 100000
$ws.Range("J46").Value = 47632.715
$ws.Range("K46").Value = 100000
$ws.Range("L46").Value = 47632.715
$ws.Range("M46").Value = -99769
$ws.Range("N46").Value = -48094.715
$ws.Range("H132").Value = 490.07407
$ws.Range("I132").Value = 403.0408
$ws.Range("J132").Value = 1343
$ws.Range("K132").Value = 1209.1224
$ws.Range("L132").Value = 4029
$ws.Range("M132").Value = 1320.8776
$ws.Range("N132").Value = -9089
$ws.Range("H134").Value = 54178.625
$ws.Range("I134").Value = 100000
$ws.Range("J134").Value = 47632.715
$ws.Range("K134").Value = 300000
$ws.Range("L134").Value = 142898.145
$ws.Range("M134").Value = -297465
$ws.Range("N134").Value = -147968.145
$ws.Range("H136").Value = 376.85
$ws.Range("I136").Value = 306.2353
$ws.Range("J136").Value = 777
$ws.Range("K136").Value = 918.7058999999999
$ws.Range("L136").Value = 2331
$ws.Range("M136").Value = 1631.2941
$ws.Range("N136").Value = -7431

Write-Output "Applied Yojimbo_Profits scheduled update."